$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row with "ChatVRM" in cell A15
$ws.Range("A15").Value = "ChatVRM"

# Update the selection to match the new active cell (A15)
$ws.Range("A15").Select()
